$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: "Name" -> "Name`nMultiline" in C1, reformat like B1 (wrap-capable header style) ---
$ws.Range("C1").Value = "Name" + [char]10 + "Multiline"
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Rows.Item(1).RowHeight = 15.75

# --- Clear the stale roster rows (26-32): this "un-joins" the extra people that were
#     previously appended, while keeping the row/column formatting intact ---
$ws.Range("A26:D32").ClearContents()

# --- Extend the formatted (but now empty) block down through row 39, matching the
#     existing style pattern: columns A, B, D use the plain row style, column C uses
#     the name-column style ---
$ws.Range("A2").Copy()
$ws.Range("A26:B39").PasteSpecial(-4122)
$ws.Range("D26:D39").PasteSpecial(-4122)

$ws.Range("B3").Copy()
$ws.Range("C26:C39").PasteSpecial(-4122)

# --- Extend the trailing placeholder rows from 1000 down to 1007 ---
for ($r = 1001; $r -le 1007; $r++) {
    $ws.Rows.Item($r).RowHeight = 15.75
}
